# Update "Remaining" (column E) values on the active worksheet to reflect
# the latest usage figures for the Slack-integration ISA report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 153.09
    3  = 139.64
    4  = 343.35
    8  = 215.74
    9  = 136.27
    10 = 52.98
    11 = 112.47
    12 = 139.25
    13 = 319.38
    15 = 383.15
    16 = 110.07
    17 = 124.23
    18 = 166.3
    19 = 135.66
    20 = 139.18
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
